# Update "Estado de Cuenta" worksheet: reverse the Periodo Mora column order
# (2008..1907 -> 1907..2008) and swap the two Valor Mora values that moved
# with the first/last period (14667 and 110000).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column E (Periodo Mora) currently runs 2008 down to 1907 for rows 16-29.
# After the edit it should run 1907 up to 2008 (fully reversed order).
$newPeriods = @("1907", "1908", "1909", "1910", "1911", "1912", "2001", "2002", "2003", "2004", "2005", "2006", "2007", "2008")

for ($i = 0; $i -lt $newPeriods.Length; $i++) {
    $row = 16 + $i
    $ws.Range("E$row").Value = $newPeriods[$i]
}

# The "Valor Mora" figures for the first and last period (rows 16 and 29)
# swap along with the reordering.
$ws.Range("F16").Value = 110000
$ws.Range("F29").Value = 14667
